$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first two data rows (rows 2 and 3), shifting all remaining
# rows up by two. This moves the starting username from adnewuser_143 to
# adnewuser_145, and the last two rows (previously 53 and 54) are dropped.
$ws.Range("A2:B3").EntireRow.Delete()
